$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 2008.875
$ws_ALC.Range("I33").Value = 2337.5
$ws_ALC.Range("K33").Value = 2337.5
$ws_ALC.Range("M33").Value = -2108.5
$ws_ALC.Range("H40").Value = 4884.85
$ws_ALC.Range("J40").Value = 4999
$ws_ALC.Range("L40").Value = 4999
$ws_ALC.Range("N40").Value = -5349
$ws_ALC.Range("H46").Value = 846.3333
$ws_ALC.Range("I46").Value = 849.5
$ws_ALC.Range("K46").Value = 2548.5
$ws_ALC.Range("M46").Value = -2429.5
$ws_ALC.Range("H60").Value = 846.3333
$ws_ALC.Range("I60").Value = 849.5
$ws_ALC.Range("K60").Value = 2548.5
$ws_ALC.Range("M60").Value = -2064.5
$ws_ALC.Range("H69").Value = 12784.8
$ws_ALC.Range("I69").Value = 12123.75
$ws_ALC.Range("J69").Value = 13225.5
$ws_ALC.Range("K69").Value = 36371.25
$ws_ALC.Range("L69").Value = 39676.5
$ws_ALC.Range("M69").Value = -35497.25
$ws_ALC.Range("N69").Value = -41424.5
$ws_ALC.Range("H72").Value = 12784.8
$ws_ALC.Range("I72").Value = 12123.75
$ws_ALC.Range("J72").Value = 13225.5
$ws_ALC.Range("K72").Value = 109113.75
$ws_ALC.Range("L72").Value = 119029.5
$ws_ALC.Range("M72").Value = -104745.75
$ws_ALC.Range("N72").Value = -127765.5
$ws_ALC.Range("H100").Value = 1758
$ws_ALC.Range("I100").Value = 1554.1666
$ws_ALC.Range("K100").Value = 1554.1666
$ws_ALC.Range("M100").Value = -1013.1666
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 4144.3276
$ws_ARM.Range("I61").Value = 4139.724
$ws_ARM.Range("J61").Value = 4233.3335
$ws_ARM.Range("K61").Value = 4139.724
$ws_ARM.Range("L61").Value = 4233.3335
$ws_ARM.Range("M61").Value = -3927.724
$ws_ARM.Range("N61").Value = -4657.3335
$ws_ARM.Range("H136").Value = 4144.3276
$ws_ARM.Range("I136").Value = 4139.724
$ws_ARM.Range("J136").Value = 4233.3335
$ws_ARM.Range("K136").Value = 12419.172
$ws_ARM.Range("L136").Value = 12700.0005
$ws_ARM.Range("M136").Value = -9869.172
$ws_ARM.Range("N136").Value = -17800.0005
$ws_ARM.Range("H137").Value = 161249.25
$ws_ARM.Range("I137").Value = 45000
$ws_ARM.Range("J137").Value = 199999
$ws_ARM.Range("K137").Value = 45000
$ws_ARM.Range("L137").Value = 199999
$ws_ARM.Range("M137").Value = -39900
$ws_ARM.Range("N137").Value = -210199
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H22").Value = 55875.11
$ws_BSM.Range("J22").Value = 167081.33
$ws_BSM.Range("L22").Value = 167081.33
$ws_BSM.Range("N22").Value = -167427.33
$ws_BSM.Range("H86").Value = 2309.4375
$ws_BSM.Range("I86").Value = 2421.5386
$ws_BSM.Range("J86").Value = 1823.6666
$ws_BSM.Range("K86").Value = 2421.5386
$ws_BSM.Range("L86").Value = 1823.6666
$ws_BSM.Range("M86").Value = -1298.5386
$ws_BSM.Range("N86").Value = -4069.6666
$ws_BSM.Range("H89").Value = 2309.4375
$ws_BSM.Range("I89").Value = 2421.5386
$ws_BSM.Range("J89").Value = 1823.6666
$ws_BSM.Range("K89").Value = 12107.693
$ws_BSM.Range("L89").Value = 9118.333000000001
$ws_BSM.Range("M89").Value = -6491.692999999999
$ws_BSM.Range("N89").Value = -20350.333
$ws_BSM.Range("H134").Value = 4686.375
$ws_BSM.Range("I134").Value = 3748.7727
$ws_BSM.Range("K134").Value = 11246.3181
$ws_BSM.Range("M134").Value = -8711.3181
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H134").Value = 5682.875
$ws_CRP.Range("I134").Value = 5597.4546
$ws_CRP.Range("K134").Value = 16792.3638
$ws_CRP.Range("M134").Value = -14257.3638
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 526.6429000000001
$ws_CUL.Range("I5").Value = 470.52
$ws_CUL.Range("K5").Value = 1411.56
$ws_CUL.Range("M5").Value = -1299.56
$ws_CUL.Range("H113").Value = 1432.826
$ws_CUL.Range("I113").Value = 1239.6666
$ws_CUL.Range("J113").Value = 1501
$ws_CUL.Range("K113").Value = 3718.9998
$ws_CUL.Range("L113").Value = 4503
$ws_CUL.Range("M113").Value = -1548.9998
$ws_CUL.Range("N113").Value = -8843
$ws_CUL.Range("H126").Value = 5000
$ws_CUL.Range("J126").Value = 5000
$ws_CUL.Range("L126").Value = 15000
$ws_CUL.Range("N126").Value = -24880
$ws_CUL.Range("H129").Value = 1646.0714
$ws_CUL.Range("I129").Value = 580
$ws_CUL.Range("K129").Value = 1740
$ws_CUL.Range("M129").Value = 3260
$ws_CUL.Range("H135").Value = 526.6429000000001
$ws_CUL.Range("I135").Value = 470.52
$ws_CUL.Range("K135").Value = 4234.68
$ws_CUL.Range("M135").Value = -1699.68
$ws_CUL.Range("H136").Value = 12571.333
$ws_CUL.Range("I136").Value = 11085.6
$ws_CUL.Range("K136").Value = 33256.8
$ws_CUL.Range("M136").Value = -28156.8
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 43209920
$ws_GSM.Range("I2").Value = 59829090
$ws_GSM.Range("K2").Value = 59829090
$ws_GSM.Range("M2").Value = -59828977
$ws_GSM.Range("H41").Value = 26617
$ws_GSM.Range("I41").Value = 29925.5
$ws_GSM.Range("K41").Value = 29925.5
$ws_GSM.Range("M41").Value = -29570.5
$ws_GSM.Range("H80").Value = 19169568
$ws_GSM.Range("I80").Value = 23002798
$ws_GSM.Range("J80").Value = 3420.5
$ws_GSM.Range("K80").Value = 23002798
$ws_GSM.Range("L80").Value = 3420.5
$ws_GSM.Range("M80").Value = -23001800
$ws_GSM.Range("N80").Value = -5416.5
$ws_GSM.Range("H83").Value = 19169568
$ws_GSM.Range("I83").Value = 23002798
$ws_GSM.Range("J83").Value = 3420.5
$ws_GSM.Range("K83").Value = 115013990
$ws_GSM.Range("L83").Value = 17102.5
$ws_GSM.Range("M83").Value = -115008998
$ws_GSM.Range("N83").Value = -27086.5
$ws_GSM.Range("H113").Value = 2400.9285
$ws_GSM.Range("I113").Value = 1419.909
$ws_GSM.Range("J113").Value = 5998
$ws_GSM.Range("K113").Value = 1419.909
$ws_GSM.Range("L113").Value = 5998
$ws_GSM.Range("M113").Value = 750.0909999999999
$ws_GSM.Range("N113").Value = -10338
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 11049.8
$ws_LTW.Range("I7").Value = 7061.875
$ws_LTW.Range("K7").Value = 7061.875
$ws_LTW.Range("M7").Value = -6949.875
$ws_LTW.Range("H55").Value = 345.08334
$ws_LTW.Range("I55").Value = 103.818184
$ws_LTW.Range("J55").Value = 2999
$ws_LTW.Range("K55").Value = 103.818184
$ws_LTW.Range("L55").Value = 2999
$ws_LTW.Range("M55").Value = 69.181816
$ws_LTW.Range("N55").Value = -3345
$ws_LTW.Range("H68").Value = 3087.5312
$ws_LTW.Range("I68").Value = 1242.8928
$ws_LTW.Range("J68").Value = 16000
$ws_LTW.Range("K68").Value = 1242.8928
$ws_LTW.Range("L68").Value = 16000
$ws_LTW.Range("M68").Value = -493.8928000000001
$ws_LTW.Range("N68").Value = -17498
$ws_LTW.Range("H71").Value = 3087.5312
$ws_LTW.Range("I71").Value = 1242.8928
$ws_LTW.Range("J71").Value = 16000
$ws_LTW.Range("K71").Value = 6214.464
$ws_LTW.Range("L71").Value = 80000
$ws_LTW.Range("M71").Value = -2470.464
$ws_LTW.Range("N71").Value = -87488
$ws_LTW.Range("H93").Value = 1224
$ws_LTW.Range("J93").Value = 0
$ws_LTW.Range("L93").Value = 0
$ws_LTW.Range("N93").ClearContents()
$ws_LTW.Range("H100").Value = 2308.875
$ws_LTW.Range("J100").Value = 2632.6667
$ws_LTW.Range("L100").Value = 2632.6667
$ws_LTW.Range("N100").Value = -3714.6667
$ws_LTW.Range("H126").Value = 11049.8
$ws_LTW.Range("I126").Value = 7061.875
$ws_LTW.Range("K126").Value = 21185.625
$ws_LTW.Range("M126").Value = -18715.625
$ws_LTW.Range("H136").Value = 6017447.5
$ws_LTW.Range("I136").Value = 10607479
$ws_LTW.Range("K136").Value = 31822437
$ws_LTW.Range("M136").Value = -31819887
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H96").Value = 93203.63
$ws_WVR.Range("I96").Value = 145491.86
$ws_WVR.Range("K96").Value = 145491.86
$ws_WVR.Range("M96").Value = -144118.86
$ws_WVR.Range("H113").Value = 808.5
$ws_WVR.Range("J113").Value = 1408.8334
$ws_WVR.Range("L113").Value = 4226.5002
$ws_WVR.Range("N113").Value = -8566.5002
